$d = $word.ActiveDocument

# "source sink rate" table: columns are
#   1: biome preference, 2: away, 3: into, 4: source sink ratio
# Rows 2-6 (row 1 is the header) get their "away"/"into"/"source sink ratio"
# values revised per the updated analysis.
$tbl = $d.Tables.Item(1)

$updates = @(
    @(2, 2, "32.70"), @(2, 3, "6.37"),  @(2, 4, "5.13"),
    @(3, 2, "28.79"), @(3, 3, "25.69"), @(3, 4, "1.12"),
    @(4, 2, "17.52"), @(4, 3, "25.38"), @(4, 4, "0.69"),
    @(5, 2, "5.66"),  @(5, 3, "23.02"), @(5, 4, "0.25"),
    @(6, 2, "0.29"),  @(6, 3, "4.50"),  @(6, 4, "0.06")
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $newValue = $u[2]
    $tbl.Cell($row, $col).Range.Text = $newValue
}
